# Add a "Label" column (H) that encodes diagnosis group as 0 (Control) / 1 (MDD)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Label"
# Match the bold/border/centered header formatting used by the other header cells (B1:G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Rows 2-11: first block (Control rows 2-6 -> 0, MDD rows 7-11 -> 1)
# Rows 12-21: second block (Control rows 12-16 -> 0, MDD rows 17-21 -> 1)
for ($r = 2; $r -le 21; $r++) {
    if (($r -ge 2 -and $r -le 6) -or ($r -ge 12 -and $r -le 16)) {
        $ws.Cells.Item($r, 8).Value = 0
    } else {
        $ws.Cells.Item($r, 8).Value = 1
    }
}
